$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.160.56"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.680.09"
$ws.Range("E3").Value = "  -3.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "681.93"
$ws.Range("E5").Value = "  -3.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.38"
$ws.Range("E6").Value = "  -4.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.679.96"
$ws.Range("E7").Value = "  -3.01%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -4.56%  "

$ws.Range("E10").Value = "  -7.53%  "

$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("E13").Value = "  -7.23%  "

$ws.Range("E14").Value = "  -7.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.301.65"
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.669.51"
$ws.Range("E16").Value = "  -4.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.215.96"
$ws.Range("E17").Value = "  -2.08%  "

$ws.Range("E18").Value = "  -1.77%  "

$ws.Range("E19").Value = "  -6.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.64"
$ws.Range("E20").Value = "  -7.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.97"
$ws.Range("E21").Value = "  -1.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.82"
$ws.Range("E22").Value = "  -7.51%  "

$ws.Range("E23").Value = "  -8.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.70"
$ws.Range("E24").Value = "  -6.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.824.30"
$ws.Range("E25").Value = "  -3.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.54"
$ws.Range("E26").Value = "  -4.44%  "

$ws.Range("E27").Value = "  -12.14%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("E29").Value = "  -8.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.85"
$ws.Range("E30").Value = "  -9.90%  "

$ws.Range("E31").Value = "  -11.29%  "

$ws.Range("E32").Value = "  -5.00%  "

$ws.Range("E33").Value = "  -8.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.94"
$ws.Range("E35").Value = "  -7.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.162"
$ws.Range("E36").Value = "  -6.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.647.05"

$ws.Range("E38").Value = "  -5.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.07"
$ws.Range("E39").Value = "  +2.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0946"
$ws.Range("E40").Value = "  -6.54%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.18"
$ws.Range("E42").Value = "  -5.65%  "

$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.957"
$ws.Range("E44").Value = "  -7.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "156.79"
$ws.Range("E45").Value = "  -4.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.97"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("E47").Value = "  -15.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "390.84"
$ws.Range("E50").Value = "  -7.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "28.15"
$ws.Range("E51").Value = "  +1.75%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.31"
$ws.Range("E48").Value = "  -3.90%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000277"
$ws.Range("E49").Value = "  -12.83%  "
